# KIBON-2480: Add BFS Nr. Standortgemeinde to Institutionen Statistik
#
# The "Institutionen" report sheet gains two new columns next to the
# existing Standortgemeinde / Traegergemeinde columns:
#   - BFS Nr. Standortgemeinde (right after Standortgemeinde)
#   - BFS Nr. Traegergemeinde  (right after Traegergemeinde)
# The old, single "BFS Gemeinde" column is replaced by these two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout before this edit (row 4 = header titles, row 5 = placeholders):
#   N = Standortgemeinde, O = Traegergemeinde, P = BFS Gemeinde, Q = Telefon ...
#
# Inserting a new column at O shifts Traegergemeinde -> P, BFS Gemeinde -> Q,
# Telefon -> R, etc. The new O becomes the BFS Nr. Standortgemeinde column,
# and the old "BFS Gemeinde" column (now sitting at Q) is overwritten in place
# to become the new BFS Nr. Traegergemeinde column - giving the final layout:
#   N = Standortgemeinde, O = BFS Nr. Standortgemeinde,
#   P = Traegergemeinde, Q = BFS Nr. Traegergemeinde, R = Telefon ...
$ws.Columns("O:O").Insert() | Out-Null

# Fill in the new / repurposed header cells (row 4) and placeholder cells
# (row 5). Write order matches the order the new placeholder strings should
# be appended to the shared string table.
$ws.Range("O4").Value = "{bfsStandortgemeindeTitle}"
$ws.Range("Q4").Value = "{bfsTraegergemeindeTitle}"
$ws.Range("Q5").Value = "{bfsTraegergemeinde}"
$ws.Range("O5").Value = "{bfsStandortgemeinde}"

# Match the reviewer's final cursor position/selection on the sheet.
$ws.Range("M11").Select() | Out-Null
